$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8 ----
$ws.Range("A8").Value = ""
$ws.Range("C8").Value = 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.0"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it's ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.00"
$ws.Range("G8").Style = "Normal"

# ---- Row 9 ----
$ws.Range("A9").Value = ""
$ws.Range("C9").Value = 65
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "11.0"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F9").Value = 0
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.00"
$ws.Range("G9").Style = "Normal"

# ---- Row 10 ----
$ws.Range("A10").Value = "Set"
$ws.Range("C10").Value = 75
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "13.0"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Range("F10").Value = 5733
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "429975.00"
$ws.Range("G10").Style = "Normal"

# ---- Row 11 ----
$ws.Range("A11").Value = ""
$ws.Range("C11").Value = 39
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.0"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F11").Value = 0
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.00"
$ws.Range("G11").Style = "Normal"

# ---- Row 12 ----
$ws.Range("A12").Value = "%"
$ws.Range("C12").Value = 78
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "Add Tender Premium "
# F12 and G12 are unchanged by the diff

# ---- Row 14 ----
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "429975.00"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "429975.00"
$ws.Range("H14").Style = "Normal"

# ---- Row 16 ----
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "429975.00"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "429975.00"
$ws.Range("H16").Style = "Normal"
